$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Photometric-Opt")
$ws1.Range("D4").Value = 7.44685327999801
$ws1.Range("K4").Value = 1.638660234711782
$ws1.Range("D5").Value = 7.402409024285747
$ws1.Range("K5").Value = 1.630821204415821
$ws1.Range("B6").Value = 7.44685327999801
$ws1.Range("C6").Value = 7.402409024285747
$ws1.Range("E6").Value = 7.508599766439922
$ws1.Range("F6").Value = 7.43725906158682
$ws1.Range("I6").Value = 1.638660234711782
$ws1.Range("J6").Value = 1.630821204415821
$ws1.Range("L6").Value = 5.150636344937815
$ws1.Range("M6").Value = 1.657017485361074
$ws1.Range("D7").Value = 7.508599766439922
$ws1.Range("K7").Value = 5.150636344937815
$ws1.Range("D8").Value = 7.43725906158682
$ws1.Range("K8").Value = 1.657017485361074
$ws1.Range("D14").Value = 7.446853279991807
$ws1.Range("K14").Value = 1.251583407461542
$ws1.Range("D15").Value = 7.402409024286854
$ws1.Range("K15").Value = 1.251648196523142
$ws1.Range("B16").Value = 7.446853279991807
$ws1.Range("C16").Value = 7.402409024286854
$ws1.Range("E16").Value = 7.508599766440905
$ws1.Range("F16").Value = 7.43725906158266
$ws1.Range("I16").Value = 1.251583407461542
$ws1.Range("J16").Value = 1.251648196523142
$ws1.Range("L16").Value = 1.249225073247281
$ws1.Range("M16").Value = 1.250475563523044
$ws1.Range("D17").Value = 7.508599766440905
$ws1.Range("K17").Value = 1.249225073247281
$ws1.Range("D18").Value = 7.43725906158266
$ws1.Range("K18").Value = 1.250475563523044

$ws2 = $wb.Worksheets.Item("Photometric-Pess")
$ws2.Range("D4").Value = 3.384266125046481
$ws2.Range("K4").Value = 1.189449767637524
$ws2.Range("D5").Value = 3.918593313466014
$ws2.Range("K5").Value = 1.202402704743217
$ws2.Range("B6").Value = 3.384266125046481
$ws2.Range("C6").Value = 3.918593313466014
$ws2.Range("E6").Value = 4.242445508897008
$ws2.Range("F6").Value = 3.402226655773111
$ws2.Range("I6").Value = 1.189449767637524
$ws2.Range("J6").Value = 1.202402704743217
$ws2.Range("L6").Value = 1.636833993271819
$ws2.Range("M6").Value = 1.269127300816686
$ws2.Range("D7").Value = 4.242445508897008
$ws2.Range("K7").Value = 1.636833993271819
$ws2.Range("D8").Value = 3.402226655773111
$ws2.Range("K8").Value = 1.269127300816686
$ws2.Range("D14").Value = 3.38426612505494
$ws2.Range("K14").Value = 1.592455224818253
$ws2.Range("D15").Value = 3.62054734826109
$ws2.Range("K15").Value = 4.097090025903984
$ws2.Range("B16").Value = 3.38426612505494
$ws2.Range("C16").Value = 3.62054734826109
$ws2.Range("E16").Value = 4.242445508900473
$ws2.Range("F16").Value = 3.402226655778946
$ws2.Range("I16").Value = 1.592455224818253
$ws2.Range("J16").Value = 4.097090025903984
$ws2.Range("L16").Value = 4.09713449033247
$ws2.Range("M16").Value = 1.597460369086457
$ws2.Range("D17").Value = 4.242445508900473
$ws2.Range("K17").Value = 4.09713449033247
$ws2.Range("D18").Value = 3.402226655778946
$ws2.Range("K18").Value = 1.597460369086457
